$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137 (ALC)
$ws.Range("H137").Value = 4168971.8
$ws.Range("I137").Value = 1357.28
$ws.Range("J137").Value = 8698988
$ws.Range("K137").Value = 4071.84
$ws.Range("L137").Value = 26096964
$ws.Range("M137").Value = -1521.84
$ws.Range("N137").Value = -26102064

# Row 138 (ALC)
$ws.Range("H138").Value = 5955276.5
$ws.Range("I138").Value = 1437.28
$ws.Range("J138").Value = 14710923
$ws.Range("K138").Value = 4311.84
$ws.Range("L138").Value = 44132769
$ws.Range("M138").Value = 828.1599999999999
$ws.Range("N138").Value = -44143049

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 8044.9663
$ws.Range("I32").Value = 9932.486000000001
$ws.Range("J32").Value = 4870.5
$ws.Range("K32").Value = 9932.486000000001
$ws.Range("L32").Value = 4870.5
$ws.Range("M32").Value = -9645.486000000001
$ws.Range("N32").Value = -5444.5

# Row 61 (ARM)
$ws.Range("H61").Value = 11907059
$ws.Range("I61").Value = 15153713
$ws.Range("J61").Value = 2660.2222
$ws.Range("K61").Value = 15153713
$ws.Range("L61").Value = 2660.2222
$ws.Range("M61").Value = -15153501
$ws.Range("N61").Value = -3084.2222

# Row 74 (ARM)
$ws.Range("H74").Value = 18520760
$ws.Range("I74").Value = 25001718
$ws.Range("J74").Value = 3740.8572
$ws.Range("K74").Value = 25001718
$ws.Range("L74").Value = 3740.8572
$ws.Range("M74").Value = -25000844
$ws.Range("N74").Value = -5488.8572

# Row 77 (ARM)
$ws.Range("H77").Value = 18520760
$ws.Range("I77").Value = 25001718
$ws.Range("J77").Value = 3740.8572
$ws.Range("K77").Value = 125008590
$ws.Range("L77").Value = 18704.286
$ws.Range("M77").Value = -125004222
$ws.Range("N77").Value = -27440.286

# Row 132 (ARM)
$ws.Range("H132").Value = 15628246
$ws.Range("I132").Value = 17859780
$ws.Range("J132").Value = 7506.5
$ws.Range("K132").Value = 53579340
$ws.Range("L132").Value = 22519.5
$ws.Range("M132").Value = -53576810
$ws.Range("N132").Value = -27579.5

# Row 136 (ARM)
$ws.Range("H136").Value = 11907059
$ws.Range("I136").Value = 15153713
$ws.Range("J136").Value = 2660.2222
$ws.Range("K136").Value = 45461139
$ws.Range("L136").Value = 7980.6666
$ws.Range("M136").Value = -45458589
$ws.Range("N136").Value = -13080.6666

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (BSM)
$ws.Range("H94").Value = 797.4706
$ws.Range("I94").Value = 735.8
$ws.Range("K94").Value = 735.8
$ws.Range("M94").Value = -284.8

# Row 99 (BSM)
$ws.Range("H99").Value = 1091.16
$ws.Range("I99").Value = 874
$ws.Range("J99").Value = 1959.8
$ws.Range("K99").Value = 874
$ws.Range("L99").Value = 1959.8
$ws.Range("M99").Value = 624
$ws.Range("N99").Value = -4955.8

# Row 113 (BSM)
$ws.Range("H113").Value = 6933.2
$ws.Range("I113").Value = 6933.2
$ws.Range("K113").Value = 6933.2
$ws.Range("M113").Value = -4763.2

# Row 134 (BSM)
$ws.Range("H134").Value = 3421.6924
$ws.Range("I134").Value = 2231.0667
$ws.Range("J134").Value = 7390.4443
$ws.Range("K134").Value = 6693.2001
$ws.Range("L134").Value = 22171.3329
$ws.Range("M134").Value = -4158.2001
$ws.Range("N134").Value = -27241.3329

$ws = $wb.Worksheets.Item("CRP")
# Row 7 (CRP)
$ws.Range("H7").Value = 143.88889
$ws.Range("I7").Value = 88.70587999999999
$ws.Range("J7").Value = 237.7
$ws.Range("K7").Value = 88.70587999999999
$ws.Range("L7").Value = 237.7
$ws.Range("M7").Value = 24.29412000000001
$ws.Range("N7").Value = -463.7

# Row 22 (CRP)
$ws.Range("H22").Value = 871.5
$ws.Range("I22").Value = 576.3333
$ws.Range("K22").Value = 576.3333
$ws.Range("M22").Value = -226.3333

# Row 31 (CRP)
$ws.Range("H31").Value = 14172948
$ws.Range("I31").Value = 14629.1
$ws.Range("J31").Value = 23611828
$ws.Range("K31").Value = 14629.1
$ws.Range("L31").Value = 23611828
$ws.Range("M31").Value = -14334.1
$ws.Range("N31").Value = -23612418

# Row 34 (CRP)
$ws.Range("H34").Value = 14172948
$ws.Range("I34").Value = 14629.1
$ws.Range("J34").Value = 23611828
$ws.Range("K34").Value = 14629.1
$ws.Range("L34").Value = 23611828
$ws.Range("M34").Value = -14427.1
$ws.Range("N34").Value = -23612232

# Row 107 (CRP)
$ws.Range("H107").Value = 336.31818
$ws.Range("I107").Value = 315.16666
$ws.Range("J107").Value = 431.5
$ws.Range("K107").Value = 315.16666
$ws.Range("L107").Value = 431.5
$ws.Range("M107").Value = 1604.83334
$ws.Range("N107").Value = -4271.5

# Row 132 (CRP)
$ws.Range("H132").Value = 3806.5
$ws.Range("I132").Value = 3836
$ws.Range("J132").Value = 3777
$ws.Range("K132").Value = 11508
$ws.Range("L132").Value = 11331
$ws.Range("M132").Value = -8978
$ws.Range("N132").Value = -16391

# Row 134 (CRP)
$ws.Range("H134").Value = 918060.7
$ws.Range("I134").Value = 2151.2632
$ws.Range("J134").Value = 3404100.5
$ws.Range("K134").Value = 6453.7896
$ws.Range("L134").Value = 10212301.5
$ws.Range("M134").Value = -3918.7896
$ws.Range("N134").Value = -10217371.5

$ws = $wb.Worksheets.Item("CUL")
# Row 97 (CUL)
$ws.Range("H97").Value = 16143.286
$ws.Range("I97").Value = 26750.75
$ws.Range("K97").Value = 80252.25
$ws.Range("M97").Value = -79756.25

# Row 107 (CUL)
$ws.Range("H107").Value = 613.65717
$ws.Range("I107").Value = 281.8387
$ws.Range("J107").Value = 877.4103
$ws.Range("K107").Value = 845.5161000000001
$ws.Range("L107").Value = 2632.2309
$ws.Range("M107").Value = 1074.4839
$ws.Range("N107").Value = -6472.2309

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (GSM)
$ws.Range("H80").Value = 19610632
$ws.Range("I80").Value = 37039704
$ws.Range("J80").Value = 2924.875
$ws.Range("K80").Value = 37039704
$ws.Range("L80").Value = 2924.875
$ws.Range("M80").Value = -37038706
$ws.Range("N80").Value = -4920.875

# Row 83 (GSM)
$ws.Range("H83").Value = 19610632
$ws.Range("I83").Value = 37039704
$ws.Range("J83").Value = 2924.875
$ws.Range("K83").Value = 185198520
$ws.Range("L83").Value = 14624.375
$ws.Range("M83").Value = -185193528
$ws.Range("N83").Value = -24608.375

# Row 97 (GSM)
$ws.Range("H97").Value = 756.6667
$ws.Range("I97").Value = 866.9
$ws.Range("J97").Value = 205.5
$ws.Range("K97").Value = 866.9
$ws.Range("L97").Value = 205.5
$ws.Range("M97").Value = -370.9
$ws.Range("N97").Value = -1197.5

# Row 132 (GSM)
$ws.Range("H132").Value = 6135.6665
$ws.Range("I132").Value = 4961.5
$ws.Range("J132").Value = 7603.375
$ws.Range("K132").Value = 14884.5
$ws.Range("L132").Value = 22810.125
$ws.Range("M132").Value = -12354.5
$ws.Range("N132").Value = -27870.125

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (LTW)
$ws.Range("H46").Value = 1088.75
$ws.Range("I46").Value = 1087.1428
$ws.Range("J46").Value = 1100
$ws.Range("K46").Value = 1087.1428
$ws.Range("L46").Value = 1100
$ws.Range("M46").Value = -899.1428000000001
$ws.Range("N46").Value = -1476

# Row 82 (LTW)
$ws.Range("H82").Value = 1702.1428
$ws.Range("I82").Value = 1519.3
$ws.Range("J82").Value = 1868.3636
$ws.Range("K82").Value = 1519.3
$ws.Range("L82").Value = 1868.3636
$ws.Range("M82").Value = -1158.3
$ws.Range("N82").Value = -2590.3636

# Row 85 (LTW)
$ws.Range("H85").Value = 1702.1428
$ws.Range("I85").Value = 1519.3
$ws.Range("J85").Value = 1868.3636
$ws.Range("K85").Value = 1519.3
$ws.Range("L85").Value = 1868.3636
$ws.Range("M85").Value = -271.3
$ws.Range("N85").Value = -4364.3636

# Row 100 (LTW)
$ws.Range("H100").Value = 2000
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1459
$ws.Range("N100").ClearContents()

# Row 122 (LTW)
$ws.Range("H122").Value = 5480.5713
$ws.Range("I122").Value = 5533.3887
$ws.Range("J122").Value = 5385.5
$ws.Range("K122").Value = 16600.1661
$ws.Range("L122").Value = 16156.5
$ws.Range("M122").Value = -14150.1661
$ws.Range("N122").Value = -21056.5

# Row 132 (LTW)
$ws.Range("H132").Value = 16138271
$ws.Range("I132").Value = 4883.6924
$ws.Range("J132").Value = 27790162
$ws.Range("K132").Value = 14651.0772
$ws.Range("L132").Value = 83370486
$ws.Range("M132").Value = -12121.0772
$ws.Range("N132").Value = -83375546

# Row 136 (LTW)
$ws.Range("H136").Value = 36590732
$ws.Range("I136").Value = 55557376
$ws.Range("J136").Value = 12202.143
$ws.Range("K136").Value = 166672128
$ws.Range("L136").Value = 36606.429
$ws.Range("M136").Value = -166669578
$ws.Range("N136").Value = -41706.429

$ws = $wb.Worksheets.Item("WVR")
# Row 100 (WVR)
$ws.Range("H100").Value = 865.3333
$ws.Range("I100").Value = 417.6
$ws.Range("J100").Value = 1425
$ws.Range("K100").Value = 835.2
$ws.Range("L100").Value = 2850
$ws.Range("M100").Value = -294.2
$ws.Range("N100").Value = -3932

# Row 109 (WVR)
$ws.Range("H109").Value = 33377
$ws.Range("J109").Value = 33377
$ws.Range("L109").Value = 33377
$ws.Range("N109").Value = -36151

# Row 132 (WVR)
$ws.Range("H132").Value = 3869.4358
$ws.Range("I132").Value = 3866
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 11598
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -9068
$ws.Range("N132").Value = -17060
